# Publication release of FormatCode 1.3.0
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Metadata ---
$ws1 = $wb.Worksheets.Item(1)

# Identifier value
$ws1.Range("B3").Value = "OID:1.3.6.1.4.1.19376.1.2.7.1"

# Version value
$ws1.Range("B4").Value = "1.3.0"

# Contact rows (11 and 12)
$ws1.Range("B11").Value = "IHE (http://ihe.net)"
$ws1.Range("B12").Value = "IHE Secretary (secretary@ihe.net)"

# --- New sheet: "Include from unknown" (IPS) ---
$ws4 = $wb.Worksheets.Item(4)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count())
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Include from unknown"

# Copy layout/formatting from the "Include from DICOM" sheet (same shape)
$ws4.Range("A1:B4").Copy($newSheet.Range("A1"))

# Overwrite with the IPS-specific content
$newSheet.Range("A1").Value = "Concept"
$newSheet.Range("B1").Value = "Description"
$newSheet.Range("A2").Value = "http://hl7.org/fhir/uv/ips/StructureDefinition/Bundle-uv-ips"
$newSheet.Range("B2").Value = "IPS (International Patient Summary)"
$newSheet.Range("A3").Value = ""
$newSheet.Range("B3").Value = ""
$newSheet.Range("A4").Value = "System URI"
$newSheet.Range("B4").Value = "urn:ietf:rfc:3986"
